$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph 2 ("Simple project using react components & useSate hook.")
#    Collapse the spell-check-split runs ("useSate" wrapped in proofErr) into
#    a single run with identical visible text.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$simpleOld = "Simple project using react components & useSate hook."
$rng = $d.Content
$rng.Find.Execute($simpleOld)
$simpleStart = $rng.Start
$simpleEnd = $rng.End

# force a run-merge: briefly perturb the text, then restore it, which
# collapses the multiple runs/proofErr markers into one plain run.
$tmp = $d.Range($simpleStart, $simpleEnd)
$tmp.Text = $simpleOld + [char]1
$tmp2 = $d.Range($simpleStart, $simpleStart + $simpleOld.Length + 1)
$tmp2.Text = $simpleOld

# ---------------------------------------------------------------------------
# 2) Paragraph 2 - "Concepts: lifting state up, form submission." becomes
#    "Concepts: lifting state up, handling form submission. "
# ---------------------------------------------------------------------------
$conceptsOld = "Concepts: lifting state up, form submission."
$conceptsNew = "Concepts: lifting state up, handling form submission. "
$rng = $d.Content
$rng.Find.Execute($conceptsOld)
$cStart = $rng.Start
$cEnd = $rng.End
$rngC = $d.Range($cStart, $cEnd)
$rngC.Text = $conceptsNew

# split "Concepts: " into its own run, separate from the rest of the sentence
$part1 = "Concepts: "
$part2 = "lifting state up, handling form submission. "
$splitPos = $cStart + $part1.Length
$afterRange = $d.Range($splitPos, $splitPos + $part2.Length)
$afterText = $afterRange.Text
$afterRange.Delete()
$insertAnchor = $d.Range($splitPos, $splitPos)
$insertAnchor.InsertAfter($afterText)

Write-Output "P2 now: [$($d.Paragraphs.Item(2).Range.Text)]"

# ---------------------------------------------------------------------------
# 3) Paragraph 3 - insert new sentence about controlled/uncontrolled
#    components (in parentheses) followed by a manual line break, before the
#    existing "Main components: Expences, ExpencesForm" text.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3Start = $p3.Range.Start

$newSentence = "(controlled & uncontrolled components: when component gets value from state or props from parent component, its controlled components & when component gets value from DOM events, its uncontrolled components.)"

$insertRange = $d.Range($p3Start, $p3Start)
$insertRange.InsertAfter($newSentence + [char]11)

Write-Output "P3 now: [$($d.Paragraphs.Item(3).Range.Text)]"

# ---------------------------------------------------------------------------
# 4) Paragraph 3 - collapse the spell-check-split runs around "Expences" /
#    "ExpencesForm" into a single plain run.
# ---------------------------------------------------------------------------
$tailOld = " Expences, ExpencesForm"
$rng = $d.Content
$rng.Find.Execute($tailOld)
$tStart = $rng.Start
$tEnd = $rng.End
$tmp = $d.Range($tStart, $tEnd)
$tmp.Text = $tailOld + [char]1
$tmp2 = $d.Range($tStart, $tStart + $tailOld.Length + 1)
$tmp2.Text = $tailOld

Write-Output "P3 final: [$($d.Paragraphs.Item(3).Range.Text)]"
